# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect newer scrape values, per commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Column F values changed on sheet "展览" (rows keyed by row number)
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 11469
    3  = 10911
    6  = 985
    8  = 55
    11 = 10607
    12 = 4098
    16 = 31
    17 = 108
    18 = 420
    19 = 11098
    20 = 10858
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Column F values changed on sheet "全部类型" (note row 3 differs slightly: 10912 vs 10911)
$sheetAllTypes = $wb.Worksheets.Item("全部类型")
$allTypesUpdates = @{
    2  = 11469
    3  = 10912
    6  = 985
    8  = 55
    11 = 10607
    12 = 4098
    16 = 31
    17 = 108
    18 = 420
    19 = 11098
    20 = 10858
}
foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
